# Auto-generated Excel COM-interop script applying the Typhon_Profits.xlsx diff.
# For each affected (sheet, row): update changed numeric cells, remove cells that
# no longer exist in the target (via ClearContents, which drops the <c> entirely),
# and set any newly-introduced cells.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 100
$ws.Range("H100").Value = 2319.75
$ws.Range("I100").Value = 1661.1538
$ws.Range("J100").Value = 3542.8572
$ws.Range("K100").Value = 1661.1538
$ws.Range("L100").Value = 3542.8572
$ws.Range("M100").Value = -1120.1538
$ws.Range("N100").Value = -4624.8572

# ALC row 116
$ws.Range("H116").Value = 3206.8333
$ws.Range("I116").Value = 1996
$ws.Range("J116").Value = 3812.25
$ws.Range("K116").Value = 1996
$ws.Range("L116").Value = 3812.25
$ws.Range("M116").Value = 1446
$ws.Range("N116").Value = -10696.25

# ALC row 129
$ws.Range("H129").Value = 782.1875
$ws.Range("I129").Value = 557.8
$ws.Range("J129").Value = 884.1818
$ws.Range("K129").Value = 1673.4
$ws.Range("L129").Value = 2652.5454
$ws.Range("M129").Value = 3326.6
$ws.Range("N129").Value = -12652.5454

# ALC row 132
$ws.Range("H132").Value = 400
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 400
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 1200
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -6260

# ALC row 137
$ws.Range("H137").Value = 2225.9333
$ws.Range("I137").Value = 2596.0715
$ws.Range("K137").Value = 7788.2145
$ws.Range("M137").Value = -5238.2145

# ALC row 138
$ws.Range("H138").Value = 2534.1724
$ws.Range("I138").Value = 723.5454999999999
$ws.Range("K138").Value = 2170.6365
$ws.Range("M138").Value = 2969.3635

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 2528.1897
$ws.Range("I32").Value = 1562.711
$ws.Range("J32").Value = 5870.231
$ws.Range("K32").Value = 1562.711
$ws.Range("L32").Value = 5870.231
$ws.Range("M32").Value = -1275.711
$ws.Range("N32").Value = -6444.231

$ws = $wb.Worksheets.Item("BSM")
# BSM row 7
$ws.Range("H7").Value = 200
$ws.Range("I7").Value = 150
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 250
$ws.Range("M7").Value = -37
$ws.Range("N7").Value = -476

# BSM row 134
$ws.Range("H134").Value = 2347.262
$ws.Range("I134").Value = 2375.513
$ws.Range("K134").Value = 7126.539
$ws.Range("M134").Value = -4591.539

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws.Range("H16").Value = 1400
$ws.Range("J16").Value = 1700
$ws.Range("L16").Value = 1700
$ws.Range("N16").Value = -2274

# CRP row 94
$ws.Range("H94").Value = 5021.857
$ws.Range("I94").Value = 2025
$ws.Range("K94").Value = 2025
$ws.Range("M94").Value = -1574

# CRP row 113
$ws.Range("H113").Value = 1400
$ws.Range("J113").Value = 1700
$ws.Range("L113").Value = 1700
$ws.Range("N113").Value = -6040

# CRP row 132
$ws.Range("H132").Value = 20639.893
$ws.Range("I132").Value = 28715.945
$ws.Range("K132").Value = 86147.83499999999
$ws.Range("M132").Value = -83617.83499999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 2
$ws.Range("H2").Value = 38.2
$ws.Range("I2").Value = 29.916666
$ws.Range("J2").Value = 50.625
$ws.Range("K2").Value = 179.499996
$ws.Range("L2").Value = 303.75
$ws.Range("M2").Value = -66.49999600000001
$ws.Range("N2").Value = -529.75

# CUL row 17
$ws.Range("H17").Value = 651
$ws.Range("J17").Value = 651
$ws.Range("L17").Value = 1953
$ws.Range("N17").Value = -2291

# CUL row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# CUL row 34
$ws.Range("H34").Value = 933.1667
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 933.1667
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").Value = 2799.5001
$ws.Range("N34").Value = -2967.5001

# CUL row 39
$ws.Range("H39").Value = 2581
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2581
$ws.Range("K39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("M39").Value = 7743
$ws.Range("N39").Value = -8331

# CUL row 55
$ws.Range("H55").Value = 1550
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 300
$ws.Range("L55").Value = 9000
$ws.Range("M55").Value = -123
$ws.Range("N55").Value = -9354

# CUL row 68
$ws.Range("H68").Value = 1475
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1475
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").Value = 4425
$ws.Range("N68").Value = -6047

# CUL row 71
$ws.Range("H71").Value = 1475
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1475
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").Value = 13275
$ws.Range("N71").Value = -21387

# CUL row 98
$ws.Range("H98").Value = 565.375
$ws.Range("J98").Value = 596.5
$ws.Range("L98").Value = 1789.5
$ws.Range("N98").Value = -4785.5

# CUL row 131
$ws.Range("H131").Value = 817.35
$ws.Range("I131").Value = 407.4
$ws.Range("J131").Value = 838.92633
$ws.Range("K131").Value = 1222.2
$ws.Range("L131").Value = 2516.77899
$ws.Range("M131").Value = 3817.8
$ws.Range("N131").Value = -12596.77899

# CUL row 136
$ws.Range("H136").Value = 3365.8
$ws.Range("I136").Value = 1276.6666
$ws.Range("J136").Value = 4261.143
$ws.Range("K136").Value = 3829.9998
$ws.Range("L136").Value = 12783.429
$ws.Range("M136").Value = 1270.0002
$ws.Range("N136").Value = -22983.429

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Range("H80").Value = 3390.4
$ws.Range("I80").Value = 2893.077
$ws.Range("J80").Value = 4314
$ws.Range("K80").Value = 2893.077
$ws.Range("L80").Value = 4314
$ws.Range("M80").Value = -1895.077
$ws.Range("N80").Value = -6310

# GSM row 83
$ws.Range("H83").Value = 3390.4
$ws.Range("I83").Value = 2893.077
$ws.Range("J83").Value = 4314
$ws.Range("K83").Value = 14465.385
$ws.Range("L83").Value = 21570
$ws.Range("M83").Value = -9473.385000000002
$ws.Range("N83").Value = -31554

# GSM row 97
$ws.Range("H97").Value = 1689.2222
$ws.Range("J97").Value = 3845.8572
$ws.Range("L97").Value = 3845.8572
$ws.Range("N97").Value = -4837.8572

$ws = $wb.Worksheets.Item("LTW")
# LTW row 93
$ws.Range("H93").Value = 3106.1875
$ws.Range("I93").Value = 3030.5
$ws.Range("J93").Value = 3232.3333
$ws.Range("K93").Value = 3030.5
$ws.Range("L93").Value = 3232.3333
$ws.Range("M93").Value = -1782.5
$ws.Range("N93").Value = -5728.3333

# LTW row 132
$ws.Range("H132").Value = 671527.4399999999
$ws.Range("I132").Value = 1005499.25
$ws.Range("J132").Value = 3583.8333
$ws.Range("K132").Value = 3016497.75
$ws.Range("L132").Value = 10751.4999
$ws.Range("M132").Value = -3013967.75
$ws.Range("N132").Value = -15811.4999

$ws = $wb.Worksheets.Item("WVR")
# WVR row 54
$ws.Range("H54").Value = 15800
$ws.Range("J54").Value = 15800
$ws.Range("L54").Value = 15800
$ws.Range("N54").Value = -16840
